# Generate Report for Handback
# Update the Status cell for file "50a64f52-764a-49f7-910a-8967ae5dbeef.md" from
# "Ready for handoff" to "Handback transform failed" on the Overview, zh-cn and de-de
# sheets, and add an Error Detail (column K) explanation on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K3").Value = "Handback file name: bu3aarpu.kgy is different with handoff file name: 50a64f52-764a-49f7-910a-8967ae5dbeef.1cb9f8c1ec8062c21a20b19c47fc35c00a01e10d.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K3").Value = "Handback file name: bu3aarpu.kgy is different with handoff file name: 50a64f52-764a-49f7-910a-8967ae5dbeef.1cb9f8c1ec8062c21a20b19c47fc35c00a01e10d.de-de."
